$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 22.39000000000006
$ws.Range("G2").Value = 0.00000736321178207433
$ws.Range("H2").Value = 0.00004086912509410851
$ws.Range("K2").Value = 5.376749367832567
$ws.Range("L2").Value = "[2.814818106198742, 7.9386806294663925]"
$ws.Range("M2").Value = 0.00004453222694689174
$ws.Range("N2").Value = 0.00004453222694689174
$ws.Range("O2").Value = -1.434000250287233
$ws.Range("P2").Value = "[-2.0126319302276956, -0.8553685703467702]"
$ws.Range("Q2").Value = 0.000001564222108507707
$ws.Range("R2").Value = 0.000003128444217015414
$ws.Range("S2").Value = 13.94647844191801
$ws.Range("T2").Value = "[12.43516608906344, 15.457790794772581]"
$ws.Range("W2").Value = 5.110030030030043
$ws.Range("X2").Value = 3.048088088088094
$ws.Range("Y2").Value = 7.171971971971992

# Row 3
$ws.Range("B3").Value = 0
$ws.Range("E3").Value = 23.04000000000016
$ws.Range("G3").Value = 0.000001032300174896861
$ws.Range("H3").Value = 0.00002097352183913913
$ws.Range("K3").Value = 6.356695094409433
$ws.Range("L3").Value = "[3.6158718506458154, 9.097518338173051]"
$ws.Range("M3").Value = 0.000007124734296093393
$ws.Range("N3").Value = 0.00001424946859218679
$ws.Range("O3").Value = 0.08176317216550011
$ws.Range("P3").Value = "[-0.4465527095192705, 0.6100790538502707]"
$ws.Range("Q3").Value = 0.7609872586610922
$ws.Range("R3").Value = 0.7609872586610922
$ws.Range("S3").Value = 12.67954245978301
$ws.Range("T3").Value = "[11.015562294678187, 14.343522624887825]"
$ws.Range("W3").Value = 22.74018018018034
$ws.Range("X3").Value = 20.80288288288303
$ws.Range("Y3").Value = 24.67747747747766
